$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 123944.2
$ws.Range("I86").Value = 206027
$ws.Range("J86").Value = 820
$ws.Range("K86").Value = 206027
$ws.Range("L86").Value = 820
$ws.Range("M86").Value = -204904
$ws.Range("N86").Value = -3066
$ws.Range("H89").Value = 123944.2
$ws.Range("I89").Value = 206027
$ws.Range("J89").Value = 820
$ws.Range("K89").Value = 1030135
$ws.Range("L89").Value = 4100
$ws.Range("M89").Value = -1024519
$ws.Range("N89").Value = -15332
$ws.Range("H131").Value = 1228.7693
$ws.Range("J131").Value = 3400
$ws.Range("L131").Value = 10200
$ws.Range("N131").Value = -20280
$ws.Range("H137").Value = 1894.8182
$ws.Range("I137").Value = 1495.625
$ws.Range("J137").Value = 2122.9285
$ws.Range("K137").Value = 4486.875
$ws.Range("L137").Value = 6368.7855
$ws.Range("M137").Value = -1936.875
$ws.Range("N137").Value = -11468.7855
$ws.Range("H138").Value = 2413
$ws.Range("J138").Value = 2432.7368
$ws.Range("L138").Value = 7298.2104
$ws.Range("N138").Value = -17578.2104

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3400.5806
$ws.Range("I32").Value = 2372.1785
$ws.Range("J32").Value = 12999
$ws.Range("K32").Value = 2372.1785
$ws.Range("L32").Value = 12999
$ws.Range("M32").Value = -2085.1785
$ws.Range("N32").Value = -13573
$ws.Range("H74").Value = 2000.8462
$ws.Range("I74").Value = 1410.8334
$ws.Range("K74").Value = 1410.8334
$ws.Range("M74").Value = -536.8334
$ws.Range("H77").Value = 2000.8462
$ws.Range("I77").Value = 1410.8334
$ws.Range("K77").Value = 7054.166999999999
$ws.Range("M77").Value = -2686.166999999999
$ws.Range("H102").Value = 1938.5
$ws.Range("I102").Value = 1938.5
$ws.Range("K102").Value = 1938.5
$ws.Range("M102").Value = -316.5
$ws.Range("H132").Value = 3299.2727
$ws.Range("I132").Value = 2482.6667
$ws.Range("J132").Value = 4279.2
$ws.Range("K132").Value = 7448.000100000001
$ws.Range("L132").Value = 12837.6
$ws.Range("M132").Value = -4918.000100000001
$ws.Range("N132").Value = -17897.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1530.0769
$ws.Range("I16").Value = 989.2
$ws.Range("J16").Value = 3333
$ws.Range("K16").Value = 989.2
$ws.Range("L16").Value = 3333
$ws.Range("M16").Value = -702.2
$ws.Range("N16").Value = -3907
$ws.Range("H31").Value = 3189.3845
$ws.Range("I31").Value = 3121.8333
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 3121.8333
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -2826.8333
$ws.Range("N31").Value = -4590
$ws.Range("H34").Value = 3189.3845
$ws.Range("I34").Value = 3121.8333
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 3121.8333
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -2919.8333
$ws.Range("N34").Value = -4404
$ws.Range("H58").Value = 1062155.6
$ws.Range("I58").Value = 1554042.6
$ws.Range("J58").Value = 2706.6924
$ws.Range("K58").Value = 1554042.6
$ws.Range("L58").Value = 2706.6924
$ws.Range("M58").Value = -1553839.6
$ws.Range("N58").Value = -3112.6924
$ws.Range("H62").Value = 2566.3333
$ws.Range("I62").Value = 2566.3333
$ws.Range("K62").Value = 2566.3333
$ws.Range("M62").Value = -1942.3333
$ws.Range("H65").Value = 2566.3333
$ws.Range("I65").Value = 2566.3333
$ws.Range("K65").Value = 12831.6665
$ws.Range("M65").Value = -9711.666499999999
$ws.Range("H68").Value = 500
$ws.Range("I68").Value = 500
$ws.Range("K68").Value = 500
$ws.Range("M68").Value = 249
$ws.Range("H71").Value = 500
$ws.Range("I71").Value = 500
$ws.Range("K71").Value = 1500
$ws.Range("M71").Value = 2244
$ws.Range("H113").Value = 1530.0769
$ws.Range("I113").Value = 989.2
$ws.Range("J113").Value = 3333
$ws.Range("K113").Value = 989.2
$ws.Range("L113").Value = 3333
$ws.Range("M113").Value = 1180.8
$ws.Range("N113").Value = -7673
$ws.Range("H122").Value = 2922.0588
$ws.Range("J122").Value = 5119.875
$ws.Range("L122").Value = 15359.625
$ws.Range("N122").Value = -20259.625
$ws.Range("H132").Value = 1492.8
$ws.Range("I132").Value = 840.69446
$ws.Range("J132").Value = 3169.6428
$ws.Range("K132").Value = 2522.08338
$ws.Range("L132").Value = 9508.928400000001
$ws.Range("M132").Value = 7.916619999999966
$ws.Range("N132").Value = -14568.9284
$ws.Range("H134").Value = 1142.8125
$ws.Range("I134").Value = 1142.8125
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3428.4375
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -893.4375
$ws.Range("H136").Value = 1062155.6
$ws.Range("I136").Value = 1554042.6
$ws.Range("J136").Value = 2706.6924
$ws.Range("K136").Value = 4662127.800000001
$ws.Range("L136").Value = 8120.0772
$ws.Range("M136").Value = -4659577.800000001
$ws.Range("N136").Value = -13220.0772

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 839.75
$ws.Range("I5").Value = 700
$ws.Range("J5").Value = 979.5
$ws.Range("K5").Value = 2100
$ws.Range("L5").Value = 2938.5
$ws.Range("M5").Value = -1988
$ws.Range("N5").Value = -3162.5
$ws.Range("H26").Value = 521
$ws.Range("I26").Value = 782
$ws.Range("J26").Value = 334.57144
$ws.Range("K26").Value = 2346
$ws.Range("L26").Value = 1003.71432
$ws.Range("M26").Value = -2058
$ws.Range("N26").Value = -1579.71432
$ws.Range("H122").Value = 857.3333
$ws.Range("J122").Value = 1226.375
$ws.Range("L122").Value = 11037.375
$ws.Range("N122").Value = -15937.375
$ws.Range("H131").Value = 856.98
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 864.2653
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2592.7959
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12672.7959
$ws.Range("H132").Value = 993
$ws.Range("J132").Value = 1091.2222
$ws.Range("L132").Value = 9820.9998
$ws.Range("N132").Value = -14880.9998
$ws.Range("H135").Value = 839.75
$ws.Range("I135").Value = 700
$ws.Range("J135").Value = 979.5
$ws.Range("K135").Value = 6300
$ws.Range("L135").Value = 8815.5
$ws.Range("M135").Value = -3765
$ws.Range("N135").Value = -13885.5
$ws.Range("H137").Value = 3888.0908
$ws.Range("I137").Value = 3128.1667
$ws.Range("J137").Value = 4800
$ws.Range("K137").Value = 9384.500100000001
$ws.Range("L137").Value = 14400
$ws.Range("M137").Value = -4284.500100000001
$ws.Range("N137").Value = -24600

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3849522.2
$ws.Range("I132").Value = 12822008
$ws.Range("J132").Value = 4171.143
$ws.Range("K132").Value = 38466024
$ws.Range("L132").Value = 12513.429
$ws.Range("M132").Value = -38463494
$ws.Range("N132").Value = -17573.429
$ws.Range("H136").Value = 10753.333
$ws.Range("J136").Value = 10753.333
$ws.Range("L136").Value = 32259.999
$ws.Range("N136").Value = -37359.999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2400
$ws.Range("I22").Value = 3080
$ws.Range("J22").Value = 2022.2222
$ws.Range("K22").Value = 3080
$ws.Range("L22").Value = 2022.2222
$ws.Range("M22").Value = -2785
$ws.Range("N22").Value = -2612.2222
$ws.Range("H27").Value = 2400
$ws.Range("I27").Value = 3080
$ws.Range("J27").Value = 2022.2222
$ws.Range("K27").Value = 3080
$ws.Range("L27").Value = 2022.2222
$ws.Range("M27").Value = -2973
$ws.Range("N27").Value = -2236.2222

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 14896
$ws.Range("L104").Value = 14896
$ws.Range("N104").Value = -21884
$ws.Range("H107").Value = 1451.75
$ws.Range("I107").Value = 1001
$ws.Range("J107").Value = 1602
$ws.Range("K107").Value = 3003
$ws.Range("L107").Value = 4806
$ws.Range("M107").Value = -1083
$ws.Range("N107").Value = -8646
$ws.Range("H132").Value = 1689.2222
$ws.Range("I132").Value = 1233.8667
$ws.Range("K132").Value = 3701.6001
$ws.Range("M132").Value = -1171.6001

Write-Host "Applied all Tonberry_Profits updates"